$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 158.0527746677399

$ws.Range("A4").Value = 18814.4835
$ws.Range("B4").Value = 18508
$ws.Range("F4").Value = 6410.254
$ws.Range("G4").Value = 6378
